$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3178
$ws.Range("J17").Value = 2763.7693
$ws.Range("L17").Value = 8291.3079
$ws.Range("N17").Value = -8627.3079
$ws.Range("H28").Value = 295
$ws.Range("I28").Value = 321.1111
$ws.Range("K28").Value = 321.1111
$ws.Range("M28").Value = 163.8889
$ws.Range("H43").Value = 1626.3077
$ws.Range("I43").Value = 1730.375
$ws.Range("K43").Value = 1730.375
$ws.Range("M43").Value = -1661.375
$ws.Range("H98").Value = 1022.9737
$ws.Range("I98").Value = 929
$ws.Range("K98").Value = 929
$ws.Range("M98").Value = 569
$ws.Range("H111").Value = 11112075
$ws.Range("I111").Value = 14286554
$ws.Range("J111").Value = 1399.5
$ws.Range("K111").Value = 42859662
$ws.Range("L111").Value = 4198.5
$ws.Range("M111").Value = -42856595
$ws.Range("N111").Value = -10332.5
$ws.Range("H122").Value = 1022.9737
$ws.Range("I122").Value = 929
$ws.Range("K122").Value = 2787
$ws.Range("M122").Value = -337

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 100000000
$ws.Range("J5").Value = 100000000
$ws.Range("L5").Value = 100000000
$ws.Range("N5").Value = -100000224
$ws.Range("H32").Value = 8605.83
$ws.Range("I32").Value = 6412.23
$ws.Range("J32").Value = 23286.076
$ws.Range("K32").Value = 6412.23
$ws.Range("L32").Value = 23286.076
$ws.Range("M32").Value = -6125.23
$ws.Range("N32").Value = -23860.076
$ws.Range("H61").Value = 27095.182
$ws.Range("I61").Value = 40259.43
$ws.Range("K61").Value = 40259.43
$ws.Range("M61").Value = -40047.43
$ws.Range("H102").Value = 2146.3333
$ws.Range("I102").Value = 2146.3333
$ws.Range("K102").Value = 2146.3333
$ws.Range("M102").Value = -524.3332999999998
$ws.Range("H122").Value = 1562.6111
$ws.Range("I122").Value = 1459.2142
$ws.Range("K122").Value = 4377.642599999999
$ws.Range("M122").Value = -1927.642599999999
$ws.Range("H132").Value = 1811.0714
$ws.Range("I132").Value = 1443.6222
$ws.Range("J132").Value = 2472.48
$ws.Range("K132").Value = 4330.8666
$ws.Range("L132").Value = 7417.440000000001
$ws.Range("M132").Value = -1800.8666
$ws.Range("N132").Value = -12477.44
$ws.Range("H136").Value = 27095.182
$ws.Range("I136").Value = 40259.43
$ws.Range("K136").Value = 120778.29
$ws.Range("M136").Value = -118228.29

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 100000000
$ws.Range("J4").Value = 100000000
$ws.Range("L4").Value = 100000000
$ws.Range("N4").Value = -100000230
$ws.Range("H54").Value = 9600
$ws.Range("I54").Value = 6400
$ws.Range("K54").Value = 6400
$ws.Range("M54").Value = -5916
$ws.Range("H99").Value = 2799.5
$ws.Range("I99").Value = 2799.5
$ws.Range("K99").Value = 2799.5
$ws.Range("M99").Value = -1301.5
$ws.Range("H134").Value = 4550.533
$ws.Range("I134").Value = 4556.054
$ws.Range("J134").Value = 4525
$ws.Range("K134").Value = 13668.162
$ws.Range("L134").Value = 13575
$ws.Range("M134").Value = -11133.162
$ws.Range("N134").Value = -18645

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2139.2778
$ws.Range("I31").Value = 1803.0714
$ws.Range("J31").Value = 3316
$ws.Range("K31").Value = 1803.0714
$ws.Range("L31").Value = 3316
$ws.Range("M31").Value = -1508.0714
$ws.Range("N31").Value = -3906
$ws.Range("H34").Value = 2139.2778
$ws.Range("I34").Value = 1803.0714
$ws.Range("J34").Value = 3316
$ws.Range("K34").Value = 1803.0714
$ws.Range("L34").Value = 3316
$ws.Range("M34").Value = -1601.0714
$ws.Range("N34").Value = -3720
$ws.Range("H107").Value = 2047.3043
$ws.Range("I107").Value = 1699.375
$ws.Range("J107").Value = 2842.5715
$ws.Range("K107").Value = 1699.375
$ws.Range("L107").Value = 2842.5715
$ws.Range("M107").Value = 220.625
$ws.Range("N107").Value = -6682.5715
$ws.Range("H132").Value = 1433.7858
$ws.Range("I132").Value = 1079.65
$ws.Range("K132").Value = 3238.95
$ws.Range("M132").Value = -708.9500000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 256.66666
$ws.Range("J23").Value = 256.66666
$ws.Range("L23").Value = 769.9999799999999
$ws.Range("N23").Value = -1239.99998
$ws.Range("H37").Value = 54499.5
$ws.Range("J37").Value = 54499.5
$ws.Range("L37").Value = 163498.5
$ws.Range("N37").Value = -163722.5
$ws.Range("H56").Value = 6656.4165
$ws.Range("I56").Value = 6656.4165
$ws.Range("K56").Value = 6656.4165
$ws.Range("M56").Value = -6126.4165
$ws.Range("H131").Value = 10018922
$ws.Range("I131").Value = 83333736
$ws.Range("J131").Value = 21447.773
$ws.Range("K131").Value = 250001208
$ws.Range("L131").Value = 64343.319
$ws.Range("M131").Value = -249996168
$ws.Range("N131").Value = -74423.319
$ws.Range("H137").Value = 5513.0835
$ws.Range("I137").Value = 3813.8
$ws.Range("J137").Value = 5960.263
$ws.Range("K137").Value = 11441.4
$ws.Range("L137").Value = 17880.789
$ws.Range("M137").Value = -6341.400000000001
$ws.Range("N137").Value = -28080.789

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2399.8572
$ws.Range("I97").Value = 2509.9
$ws.Range("K97").Value = 2509.9
$ws.Range("M97").Value = -2013.9
$ws.Range("H122").Value = 1117.6666
$ws.Range("I122").Value = 715
$ws.Range("K122").Value = 2145
$ws.Range("M122").Value = 305
$ws.Range("H126").Value = 2264858
$ws.Range("I126").Value = 2780613.2
$ws.Range("K126").Value = 8341839.600000001
$ws.Range("M126").Value = -8339369.600000001
$ws.Range("H132").Value = 898392.0600000001
$ws.Range("I132").Value = 1206165.2
$ws.Range("K132").Value = 3618495.6
$ws.Range("M132").Value = -3615965.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1004.75
$ws.Range("I22").Value = 811.6
$ws.Range("J22").Value = 1197.9
$ws.Range("K22").Value = 811.6
$ws.Range("L22").Value = 1197.9
$ws.Range("M22").Value = -516.6
$ws.Range("N22").Value = -1787.9
$ws.Range("H27").Value = 1004.75
$ws.Range("I27").Value = 811.6
$ws.Range("J27").Value = 1197.9
$ws.Range("K27").Value = 811.6
$ws.Range("L27").Value = 1197.9
$ws.Range("M27").Value = -704.6
$ws.Range("N27").Value = -1411.9
$ws.Range("H40").Value = 18536.363
$ws.Range("I40").Value = 27599.6
$ws.Range("K40").Value = 27599.6
$ws.Range("M40").Value = -27463.6
$ws.Range("H46").Value = 2203.4546
$ws.Range("I46").Value = 1798
$ws.Range("J46").Value = 2541.3333
$ws.Range("K46").Value = 1798
$ws.Range("L46").Value = 2541.3333
$ws.Range("M46").Value = -1610
$ws.Range("N46").Value = -2917.3333
$ws.Range("H68").Value = 2275.5
$ws.Range("I68").Value = 1840.6666
$ws.Range("K68").Value = 1840.6666
$ws.Range("M68").Value = -1091.6666
$ws.Range("H70").Value = 20081.5
$ws.Range("J70").Value = 20081.5
$ws.Range("L70").Value = 20081.5
$ws.Range("N70").Value = -20621.5
$ws.Range("H71").Value = 2275.5
$ws.Range("I71").Value = 1840.6666
$ws.Range("K71").Value = 9203.333000000001
$ws.Range("M71").Value = -5459.333000000001
$ws.Range("H73").Value = 20081.5
$ws.Range("J73").Value = 20081.5
$ws.Range("L73").Value = 20081.5
$ws.Range("N73").Value = -21953.5
$ws.Range("H74").Value = 47110
$ws.Range("J74").Value = 47110
$ws.Range("L74").Value = 47110
$ws.Range("N74").Value = -49106
$ws.Range("H77").Value = 47110
$ws.Range("J77").Value = 47110
$ws.Range("L77").Value = 141330
$ws.Range("N77").Value = -151314
$ws.Range("H100").Value = 1899.4286
$ws.Range("I100").Value = 1539.2
$ws.Range("K100").Value = 1539.2
$ws.Range("M100").Value = -998.2
$ws.Range("H122").Value = 2934.875
$ws.Range("I122").Value = 2746.5
$ws.Range("K122").Value = 8239.5
$ws.Range("M122").Value = -5789.5
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H62").Value = 33336666
$ws.Range("I62").Value = 33336666
$ws.Range("K62").Value = 33336666
$ws.Range("M62").Value = -33336042
$ws.Range("H65").Value = 33336666
$ws.Range("I65").Value = 33336666
$ws.Range("K65").Value = 166683330
$ws.Range("M65").Value = -166680210
$ws.Range("H70").Value = 47109
$ws.Range("J70").Value = 47109
$ws.Range("L70").Value = 47109
$ws.Range("N70").Value = -47739
$ws.Range("H73").Value = 47109
$ws.Range("J73").Value = 47109
$ws.Range("L73").Value = 47109
$ws.Range("N73").Value = -49293
$ws.Range("H132").Value = 1763.6666
$ws.Range("I132").Value = 1594.5897
$ws.Range("K132").Value = 4783.7691
$ws.Range("M132").Value = -2253.7691
